# Lesson 3 finished log update
# Updates the "Week 2" time recording sheet:
#   - Record #1 (row 7): Stop time moved from 22:30 to 23:00, and an
#     Interruption Time of 40 minutes is logged (Delta Time recalculates
#     automatically via the existing formula).
#   - Record #2 (row 8): Date logged as 06/02/2019.
#   - Selection left on B9 (next date entry cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 2")

# Row 7: Stop (D7) and Interruption Time (E7)
$ws.Range("D7").Value = 0.95833333333333337
$ws.Range("E7").Value = 40

# Row 8: Date (B8)
$ws.Range("B8").Value = 43502

# Leave the active selection on B9, as in the saved workbook
$ws.Activate()
$ws.Range("B9").Select()
